$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 424.22223
$ws.Range("I39").Value = 87.36364
$ws.Range("K39").Value = 262.09092
$ws.Range("M39").Value = 33.90908000000002

$ws.Range("H40").Value = 6527.4287
$ws.Range("I40").Value = 4423.25
$ws.Range("K40").Value = 4423.25
$ws.Range("M40").Value = -4248.25

$ws.Range("H100").Value = 2275.6667
$ws.Range("I100").Value = 2214.2917
$ws.Range("J100").Value = 2766.6667
$ws.Range("K100").Value = 2214.2917
$ws.Range("L100").Value = 2766.6667
$ws.Range("M100").Value = -1673.2917
$ws.Range("N100").Value = -3848.6667

$ws.Range("H111").Value = 2347.0908
$ws.Range("I111").Value = 2294.7144
$ws.Range("J111").Value = 2438.75
$ws.Range("K111").Value = 6884.1432
$ws.Range("L111").Value = 7316.25
$ws.Range("M111").Value = -3817.1432
$ws.Range("N111").Value = -13450.25

$ws.Range("H132").Value = 2352.5151
$ws.Range("I132").Value = 2249.7742
$ws.Range("J132").Value = 3945
$ws.Range("K132").Value = 6749.3226
$ws.Range("L132").Value = 11835
$ws.Range("M132").Value = -4219.3226
$ws.Range("N132").Value = -16895

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 58630.8
$ws.Range("J24").Value = 58630.8
$ws.Range("L24").Value = 58630.8
$ws.Range("N24").Value = -59378.8

$ws.Range("H32").Value = 4260.7046
$ws.Range("I32").Value = 2377.6829
$ws.Range("K32").Value = 2377.6829
$ws.Range("M32").Value = -2090.6829

$ws.Range("H74").Value = 2499.5
$ws.Range("I74").Value = 2499.5
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 2499.5
$ws.Range("L74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -1625.5

$ws.Range("H77").Value = 2499.5
$ws.Range("I77").Value = 2499.5
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 12497.5
$ws.Range("L77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -8129.5

$ws.Range("H92").Value = 60000
$ws.Range("I92").Value = 60000
$ws.Range("K92").Value = 60000
$ws.Range("M92").Value = -57504

$ws.Range("H100").Value = 58630.8
$ws.Range("J100").Value = 58630.8
$ws.Range("L100").Value = 58630.8
$ws.Range("N100").Value = -60794.8

$ws.Range("H107").Value = 30000
$ws.Range("J107").Value = 30000
$ws.Range("L107").Value = 30000
$ws.Range("N107").Value = -37680

$ws.Range("H122").Value = 994.5294
$ws.Range("I122").Value = 994.5294
$ws.Range("K122").Value = 2983.5882
$ws.Range("M122").Value = -533.5882000000001

$ws.Range("H132").Value = 2674.2666
$ws.Range("I132").Value = 2500.375
$ws.Range("K132").Value = 7501.125
$ws.Range("M132").Value = -4971.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 6705.875
$ws.Range("J20").Value = 14999.5
$ws.Range("L20").Value = 14999.5
$ws.Range("N20").Value = -15493.5

$ws.Range("H134").Value = 1960.4445
$ws.Range("I134").Value = 1960.4445
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 5881.333500000001
$ws.Range("L134").Value = 0
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = -3346.333500000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 7174.8
$ws.Range("I99").Value = 6589.6
$ws.Range("K99").Value = 6589.6
$ws.Range("M99").Value = -5091.6

$ws.Range("H126").Value = 7174.8
$ws.Range("I126").Value = 6589.6
$ws.Range("K126").Value = 19768.8
$ws.Range("M126").Value = -17298.8

$ws.Range("H132").Value = 3161.6843
$ws.Range("I132").Value = 1760.3572
$ws.Range("K132").Value = 5281.071599999999
$ws.Range("M132").Value = -2751.071599999999

$ws.Range("H134").Value = 1985.3077
$ws.Range("I134").Value = 1651
$ws.Range("J134").Value = 5997
$ws.Range("K134").Value = 4953
$ws.Range("L134").Value = 17991
$ws.Range("M134").Value = -2418
$ws.Range("N134").Value = -23061

$ws.Range("H141").Value = 281400.9
$ws.Range("J141").Value = 307112.12
$ws.Range("L141").Value = 307112.12
$ws.Range("N141").Value = -317472.12

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 1620.625
$ws.Range("I18").Value = 988.3333
$ws.Range("K18").Value = 2964.9999
$ws.Range("M18").Value = -2795.9999

$ws.Range("H98").Value = 631
$ws.Range("I98").Value = 549
$ws.Range("J98").Value = 672
$ws.Range("K98").Value = 1647
$ws.Range("L98").Value = 2016
$ws.Range("M98").Value = -149
$ws.Range("N98").Value = -5012

$ws.Range("H132").Value = 1904.7693
$ws.Range("I132").Value = 1319.1111
$ws.Range("J132").Value = 3222.5
$ws.Range("K132").Value = 11871.9999
$ws.Range("L132").Value = 29002.5
$ws.Range("M132").Value = -9341.999900000001
$ws.Range("N132").Value = -34062.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 19999
$ws.Range("I12").Value = 19999
$ws.Range("K12").Value = 19999
$ws.Range("M12").Value = -19859

$ws.Range("H95").Value = 13332.333
$ws.Range("J95").Value = 13332.333
$ws.Range("L95").Value = 13332.333
$ws.Range("N95").Value = -18824.333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2451.5789
$ws.Range("I46").Value = 450
$ws.Range("J46").Value = 2562.7778
$ws.Range("K46").Value = 450
$ws.Range("L46").Value = 2562.7778
$ws.Range("M46").Value = -262
$ws.Range("N46").Value = -2938.7778

$ws.Range("H55").Value = 1023.9474
$ws.Range("I55").Value = 232.25
$ws.Range("J55").Value = 1599.7273
$ws.Range("K55").Value = 232.25
$ws.Range("L55").Value = 1599.7273
$ws.Range("M55").Value = -59.25
$ws.Range("N55").Value = -1945.7273

$ws.Range("H93").Value = 1266.6666
$ws.Range("I93").Value = 900
$ws.Range("J93").Value = 2000
$ws.Range("K93").Value = 900
$ws.Range("L93").Value = 2000
$ws.Range("M93").Value = 348
$ws.Range("N93").Value = -4496

$ws.Range("H136").Value = 3283.0625
$ws.Range("I136").Value = 2902.0715
$ws.Range("K136").Value = 8706.2145
$ws.Range("M136").Value = -6156.2145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2308
$ws.Range("I126").Value = 1943.4445
$ws.Range("J126").Value = 3948.5
$ws.Range("K126").Value = 5830.333500000001
$ws.Range("L126").Value = 11845.5
$ws.Range("M126").Value = -3360.333500000001
$ws.Range("N126").Value = -16785.5

$ws.Range("H141").Value = 64537.25
$ws.Range("I141").Value = 650
$ws.Range("K141").Value = 650
$ws.Range("M141").Value = 4530
